$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (F,G) for M_TotalTax and M_CorpTax, shifting existing F:M to H:O
$ws.Range("F1:G1").EntireColumn.Insert()

# Header labels for the new columns
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# Fill in new data for M_TotalTax (F) and M_CorpTax (G) columns, rows 2-8
$ws.Range("F2").Value = 5453258922371.839
$ws.Range("G2").Value = 901477052802.8717
$ws.Range("F3").Value = 7186470855993.515
$ws.Range("G3").Value = 480200693963.0637
$ws.Range("F4").Value = 1002696890625.219
$ws.Range("G4").Value = 134030001792.1409
$ws.Range("F5").Value = 241723618867.4443
$ws.Range("G5").Value = 38152875196.83905
$ws.Range("F6").Value = 5526698445364.336
$ws.Range("G6").Value = 408160442622.5067
$ws.Range("F7").Value = 457036063703.7742
$ws.Range("G7").Value = 4932505470.985653
$ws.Range("F8").Value = 114563677189.0651
$ws.Range("G8").Value = 24178295225.23782

# Updated M_POP value for the last region (row 8)
$ws.Range("E8").Value = 366265684.25
